# Rename "Requested quantity" headers
$wb = $excel.ActiveWorkbook
$wsWeekly = $wb.Worksheets.Item(1)
$wsMonthly = $wb.Worksheets.Item(2)

$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after "Monthly Trend"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "PO Forecast"

# Match the page margins / outline settings used on the other sheets
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Copy the header style (bold, centered, bordered) from the Weekly sheet onto
# all four header cells, then overwrite the header text.
$wsWeekly.Range("A1:B1").Copy($newSheet.Range("A1:B1"))
$wsWeekly.Range("A1:B1").Copy($newSheet.Range("C1:D1"))

$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Column A (ds) rows 2-35 reuse the exact date values (and date-formatted
# style) already present on the Weekly sheet's A2:A35.
$wsWeekly.Range("A2:A35").Copy($newSheet.Range("A2:A35"))

# Rows 36-43 are new dates not present elsewhere - copy the date style from
# A2 and then overwrite with the correct values.
$wsWeekly.Range("A2").Copy($newSheet.Range("A36:A43"))

$data = @(
    @(44934.99999999999, 93, -60.73327432583762, 243.8973802666676),
    @(44941.99999999999, 95, -38.47050891951249, 236.2164013787937),
    @(44948.99999999999, 96, -44.48306436632997, 247.659437715117),
    @(44976.99999999999, 102, -39.71778977844529, 244.1720762598679),
    @(44983.99999999999, 104, -41.07518670681191, 244.3141073626756),
    @(44990.99999999999, 105, -32.85870994015621, 237.4187594373126),
    @(44997.99999999999, 107, -39.15296850439334, 260.9535779271958),
    @(45004.99999999999, 109, -39.6051090030227, 251.0405711659447),
    @(45011.99999999999, 110, -29.43012335681839, 259.275333876062),
    @(45018.99999999999, 112, -29.40090487495903, 255.9416145617073),
    @(45025.99999999999, 113, -23.72599337300892, 266.0773525194163),
    @(45060.99999999999, 121, -27.67786022349889, 254.9633444006974),
    @(45067.99999999999, 123, -23.60073445025208, 267.0101306277353),
    @(45074.99999999999, 124, -8.867742939026236, 268.9881808955302),
    @(45081.99999999999, 126, -4.696554902849408, 263.1763420932032),
    @(45088.99999999999, 127, -16.83446128956775, 275.5568977315875),
    @(45095.99999999999, 129, -18.69335252455391, 270.2455364147551),
    @(45102.99999999999, 130, -7.662187233395935, 268.2502109012177),
    @(45116.99999999999, 133, -7.253167055097542, 278.3711869233832),
    @(45123.99999999999, 135, -23.99054752529781, 273.2310966189456),
    @(45130.99999999999, 136, 7.101041968064066, 284.961731984683),
    @(45137.99999999999, 138, -6.014054751756365, 288.5886599672523),
    @(45200.99999999999, 152, 10.19146388832105, 299.8253280027792),
    @(45207.99999999999, 154, 5.546339989153595, 290.9131500324114),
    @(45214.99999999999, 155, 26.24334569549202, 295.6354168852458),
    @(45221.99999999999, 157, 19.71543439148746, 302.5577106704474),
    @(45228.99999999999, 158, 14.06001716994365, 290.7303753878352),
    @(45235.99999999999, 160, 15.80037465743688, 302.3815811487132),
    @(45249.99999999999, 163, 27.38006609313175, 309.6444686013186),
    @(45256.99999999999, 164, 24.0141453732759, 316.1677749906654),
    @(45270.99999999999, 167, 22.0942925377156, 313.7237935690571),
    @(45277.99999999999, 169, 23.22288268724595, 307.3582516713693),
    @(45298.99999999999, 174, 23.17566455413375, 317.0016968150854),
    @(45312.99999999999, 177, 40.91643597167072, 324.5660225707098),
    @(45319.99999999999, 178, 32.20337172569548, 318.2170067758544),
    @(45326.99999999999, 180, 35.47638236893621, 316.0374246658391),
    @(45333.99999999999, 181, 40.43000357561458, 326.7891263109462),
    @(45340.99999999999, 183, 43.42310891818809, 317.5399221259033),
    @(45347.99999999999, 185, 37.61501836550664, 328.7871292625424),
    @(45354.99999999999, 186, 48.66021610472106, 332.1998813019627),
    @(45361.99999999999, 188, 44.39064998305025, 330.7532238821299),
    @(45368.99999999999, 189, 45.26470752677714, 340.7421078184347)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Keep "Weekly Quantity" as the active sheet, same as before the edit.
$wsWeekly.Activate()
